$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new columns (J:N) to make room for VFX, Casting, SFX, Editing, Directing
# This shifts old J(Keywords)..Q(Genres) to O..V, and keeps I (now relabeled "Script") in place.
$ws.Range("J1:N1").EntireColumn.Insert()

# --- Header row ---
$ws.Range("I1").Value = 'Script'
$ws.Range("J1").Value = 'VFX'
$ws.Range("K1").Value = 'Casting'
$ws.Range("L1").Value = 'SFX'
$ws.Range("M1").Value = 'Editing'
$ws.Range("N1").Value = 'Directing'

# --- Row data: Summary (H), Script(I) VFX(J) Casting(K) SFX(L) Editing(M) Directing(N), Keywords(O) ---
# Row 2
$ws.Range("H2").Value = '"Everything Everywhere All at Once" is a visually stunning and emotionally resonant film that explores themes of multiverses, family dynamics, and the significance of every choice we make. With a mix of humor, action, and heartfelt moments, this movie has left a lasting impact on viewers, making it a must-watch for audiences of all ages. The film has been praised for its originality, technical brilliance, and powerful storytelling, making it a standout in the world of cinema.'
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 8
$ws.Range("K2").Value = 9
$ws.Range("L2").Value = 8
$ws.Range("M2").Value = 8
$ws.Range("N2").Value = 9
$ws.Range("O2").Value = '- audacious - funny - original - technically dazzling - thematically resonant - genetically engineered - greatest multiverse movie - required viewing - everything matters - overwhelming - in love - homage - emotional core - chaotic balance - intelligent filmmaking - lowbrow humor - hot right now - competing realities - haunting - mommy issues - nihilist lesbian representation - pride month'

# Row 3
$ws.Range("H3").Value = 'Spider-Man: Into the Spider-Verse is hailed as the best Spider-Man film, with stunning animation and a deep understanding of the character. The movie reinvigorates the superhero genre and is praised for its creativity and humor. Fans love the diverse cast of Spider-People and are excited for a sequel. The film is considered a groundbreaking achievement in animation and a must-watch for all audiences.'
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 10
$ws.Range("K3").Value = 9
$ws.Range("L3").Value = 9
$ws.Range("M3").Value = 9
$ws.Range("N3").Value = 10
$ws.Range("O3").Value = 'best - animated - creative - understanding - funny - diverse - reinvigorates - hyper-popular - sequel - overwhelming - achievement - expressive - textures - rhythms - colors - feminist - montage - greatest - favorite - paintings'

# Row 4
$ws.Range("H4").Value = 'Inception is a complex heist film dressed in science fiction conventions, following Dom Cobb as he tries to free himself from his past. It is a thought-provoking, layered story with sumptuous aesthetics and a brilliant cast, driven by Christopher Nolan''s confident direction. The film explores the idea of dreams becoming reality and leaves viewers questioning the truth behind the spinning totem at the end. Despite some bad dialogue, it is still considered a masterpiece by many.'
$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 9
$ws.Range("K4").Value = 8
$ws.Range("L4").Value = 8
$ws.Range("M4").Value = 9
$ws.Range("N4").Value = 9
$ws.Range("O4").Value = '- complex- heist- science fiction- study- man- free- past- cerebral- pop-masterpiece- thought-provoking- layered- story-telling- sumptuous- aesthetics- flawless- editing- sound design- effects- musical score- brilliant- unrivaled- filmmaking- rent free- chemistry- masterpiece- gay- lesbian- solidarity- fanfic- totem- spinning- dreaming- joyous- projection- persona 5- sexy- bad dialogue- dreams- inception'

# Row 5
$ws.Range("H5").Value = 'Spider-Man: Across the Spider-Verse is described as an overwhelming viewing experience, with viewers unable to tear their eyes away from the screen for its entire duration. The animation, humor, soundtrack, and plot are all praised as near perfect, leaving viewers beaming with joy. Some viewers were left speechless and overstimulated after the movie, with one reviewer even considering quitting making live-action films after seeing it. Overall, the film is hailed as a masterpiece, with some viewers even calling for a 6-star rating.'
$ws.Range("I5").Value = 8
$ws.Range("J5").Value = 10
$ws.Range("K5").Value = 9
$ws.Range("L5").Value = 9
$ws.Range("M5").Value = 8
$ws.Range("N5").Value = 9
$ws.Range("O5").Value = '- overwhelming - outstanding - joy - perfect - magical - mind-blowing - incredible - favorite - fried - masterpiece - ambitious - beauty - heart - funny - visually dazzling - busy - funny - psychopath - unforgettable'

# Row 6
$ws.Range("H6").Value = 'Spider-Man: No Way Home brings back Willem Dafoe''s iconic Green Goblin and delivers emotional moments with Andrew Garfield''s Spider-Man. The film is praised for its ambitious storytelling and nostalgic callbacks, but some criticize the excessive CGI and pacing issues. Overall, it is considered a thrilling and emotional experience for fans of the franchise.'
$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 8
$ws.Range("K6").Value = 9
$ws.Range("L6").Value = 8
$ws.Range("M6").Value = 6
$ws.Range("N6").Value = 5
$ws.Range("O6").Value = '- ambitious - emotional - messy - CGI - iconic - sinister - nostalgic - thrilling - surprising - ambitious - action-packed - comedic - intense - unexpected - fan-favorite - epic - tear-jerking - heartwarming - chaotic - conflicted - childhood memories - cameo - applause'

